$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.092.03'
$ws.Range("E2").Value = '  +2.37%  '
$ws.Range("D3").Value = '3.227.82'
$ws.Range("E3").Value = '  +5.65%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.59'
$ws.Range("E5").Value = '  +4.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.68'
$ws.Range("E6").Value = '  +7.26%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '3.220.27'
$ws.Range("E8").Value = '  +5.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.512'
$ws.Range("E9").Value = '  +4.49%  '
$ws.Range("E10").Value = '  +9.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.162'
$ws.Range("E11").Value = '  +5.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.488'
$ws.Range("E12").Value = '  +5.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.45'
$ws.Range("E13").Value = '  +6.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000233'
$ws.Range("E14").Value = '  +5.79%  '
$ws.Range("D15").Value = '3.753.31'
$ws.Range("E15").Value = '  +5.83%  '
$ws.Range("D16").Value = '66.182.66'
$ws.Range("E16").Value = '  +2.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '539.85'
$ws.Range("E17").Value = '  +11.08%  '
$ws.Range("D18").Value = '3.232.42'
$ws.Range("E18").Value = '  +5.71%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.114'
$ws.Range("E19").Value = '  +2.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.13'
$ws.Range("E20").Value = '  +7.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.60'
$ws.Range("E21").Value = '  +7.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.744'
$ws.Range("E22").Value = '  +8.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.77'
$ws.Range("E23").Value = '  +8.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.52'
$ws.Range("E24").Value = '  +7.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.18'
$ws.Range("E25").Value = '  +2.85%  '
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.42'
$ws.Range("E27").Value = '  +21.19%  '
$ws.Range("E28").Value = '  +8.01%  '
$ws.Range("E29").Value = '  +8.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.58'
$ws.Range("E30").Value = '  +6.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.75'
$ws.Range("E31").Value = '  +6.18%  '
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("E33").Value = '  +5.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '561.46'
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.36'
$ws.Range("E35").Value = '  +7.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.62'
$ws.Range("E36").Value = '  +4.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0456'
$ws.Range("E37").Value = '  +9.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '54.78'
$ws.Range("E38").Value = '  +4.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0861'
$ws.Range("E39").Value = '  +7.45%  '
$ws.Range("E40").Value = '  +7.60%  '
$ws.Range("D41").Value = '3.212.27'
$ws.Range("E41").Value = '  +11.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.88'
$ws.Range("E42").Value = '  +3.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.60'
$ws.Range("E43").Value = '  +4.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.284'
$ws.Range("E44").Value = '  +16.68%  '
$ws.Range("E45").Value = '  +13.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.46'
$ws.Range("E46").Value = '  +6.73%  '
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("D48").Value = '0.0₃0554'
$ws.Range("E48").Value = '  +3.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.07'
$ws.Range("E49").Value = '  +4.56%  '
$ws.Range("E50").Value = '  +3.92%  '
$ws.Range("E51").Value = '  +8.16%  '
